# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 177-178),
# pushing the existing rows 177-218 down to 179-220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 177:218 down by 2 (creates blank rows 177 and 178)
$ws.Rows("177:178").Insert()

# New row 177
$ws.Range("A177").Value = 9
$ws.Range("B177").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C177").Value = "Metropolitana"
$ws.Range("D177").Value2 = 44855
$ws.Range("E177").Value = 13
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100101
$ws.Range("H177").Value = "Berries"
$ws.Range("I177").Value = 100101001
$ws.Range("J177").Value = "Arándano (blue)"
$ws.Range("K177").Value = "Sin especificar"
$ws.Range("L177").Value = "Especial"
$ws.Range("M177").Value = 280
$ws.Range("N177").Value = 12000
$ws.Range("O177").Value = 12000
$ws.Range("P177").Value = 12000
$ws.Range("Q177").Value = "$/bandeja 2 kilos"
$ws.Range("R177").Value = "Provincia de Linares"
$ws.Range("S177").Value = 6000
$ws.Range("T177").Value = 2

# New row 178
$ws.Range("A178").Value = 9
$ws.Range("B178").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C178").Value = "Metropolitana"
$ws.Range("D178").Value2 = 44855
$ws.Range("E178").Value = 13
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100101
$ws.Range("H178").Value = "Berries"
$ws.Range("I178").Value = 100101001
$ws.Range("J178").Value = "Arándano (blue)"
$ws.Range("K178").Value = "Sin especificar"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 300
$ws.Range("N178").Value = 11000
$ws.Range("O178").Value = 11000
$ws.Range("P178").Value = 11000
$ws.Range("Q178").Value = "$/bandeja 2 kilos"
$ws.Range("R178").Value = "Provincia de Linares"
$ws.Range("S178").Value = 5500
$ws.Range("T178").Value = 2
